$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("4:5").Insert(1, 0)
